# Time warp the example model's dates from 2014 to 2021 (+7 years, same
# month/day). This mirrors the author's commit "time warp for example
# model from 2014 to 2021": a handful of date-serial cells on the
# "calendar bucket", "demand" and "parameter" sheets are bumped forward.

$wb = $excel.ActiveWorkbook

# calendar bucket: "capacity resource B" start date, "capacity resource C"
# start/end dates.
$wsCalBucket = $wb.Worksheets.Item("calendar bucket")
$wsCalBucket.Range("B3").Value = 44287
$wsCalBucket.Range("B4").Value = 44197
$wsCalBucket.Range("C4").Value = 45291

# demand: due dates for Demand 1/2/3.
$wsDemand = $wb.Worksheets.Item("demand")
$wsDemand.Range("E2").Value = 44197
$wsDemand.Range("E3").Value = 44317
$wsDemand.Range("E4").Value = 44409

# parameter: "currentdate" value.
$wsParameter = $wb.Worksheets.Item("parameter")
$wsParameter.Range("B2").Value = 44197
